# Add data for 2022-10-11: refresh the "through" date label, update the
# October figure and the running Total for 2022 in column I.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-10-11"

# Update the 2022 column header (shared string shown in I1)
$ws.Range("I1").Value = "2022 (through 10-11)"

# Update October 2022 (row 11) and the Total row (row 14) for the 2022 column
$ws.Range("I11").Value = 37
$ws.Range("I14").Value = 1315
